# Regenerate the test-data "run number" from 19/24 -> 20/25 across the
# three sheets (login, order, Sheet1) and reset the "order" sheet's
# selection/view.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet1: the two numeric "seed" cells that the CONCATENATE formulas
# in columns A/B/C (rows 2:21 and 23:42) are built from.
# ---------------------------------------------------------------
$data = $wb.Worksheets.Item("Sheet1")
$data.Range("I2").Value = 20
$data.Range("I23").Value = 25

# ---------------------------------------------------------------
# login: static (non-formula) copies of the first block of names,
# columns G/H (username) and I (email), rows 2-21, suffix 19 -> 20.
# ---------------------------------------------------------------
$login = $wb.Worksheets.Item("login")
$loginNames = @(
    "EthanBaker", "DelanieCarman", "BretAgnew", "EdgardoTaylor", "TyrekReis",
    "LeannaChow", "TuckerCarlson", "AnnmarieConnor", "MoniqueWitte", "MikelWhitlock",
    "VincentAmaya", "KeiraQuiroz", "EllisCreech", "DionteCreel", "NicholeFoust",
    "ManuelConnell", "LourdesElam", "LincolnFrederick", "AlisaCash", "LucilleGriffiths"
)
for ($idx = 0; $idx -lt $loginNames.Length; $idx++) {
    $r = $idx + 2
    $name = $loginNames[$idx] + "20"
    $email = $loginNames[$idx] + "20@gmail.com"
    $login.Cells.Item($r, 7).Value = $name
    $login.Cells.Item($r, 8).Value = $name
    $login.Cells.Item($r, 9).Value = $email
}

# ---------------------------------------------------------------
# order: static (non-formula) copies of the second block of names,
# columns R/S (username) and T (email), rows 2-21, suffix 24 -> 25.
# ---------------------------------------------------------------
$order = $wb.Worksheets.Item("order")
$orderNames = @(
    "DonnellJernigan", "MalikOtoole", "AlanCaudill", "AdanApplegate", "AiyanaWhitworth",
    "MercedezBrien", "DuaneHager", "LorenBell", "GeraldHiller", "DeionBranch",
    "DakotaHalstead", "ElliottFurman", "MiltonCamp", "DawnChester", "ZacheryPetrie",
    "EstebanAngel", "JimmyBlankenship", "AllysaGrice", "AugustineYoo", "BrandiSouthard"
)
for ($idx = 0; $idx -lt $orderNames.Length; $idx++) {
    $r = $idx + 2
    $name = $orderNames[$idx] + "25"
    $email = $orderNames[$idx] + "25@gmail.com"
    $order.Cells.Item($r, 18).Value = $name
    $order.Cells.Item($r, 19).Value = $name
    $order.Cells.Item($r, 20).Value = $email
}

# Reset the "order" sheet's active view: drop the frozen/top-left cell
# and collapse the selection down to B2.
$order.Activate()
$order.Range("B2").Select()
